# Fix the poster footer: correct the institute name spelling, and
# reposition/resize the footer text box + footer divider line so the
# longer corrected text still fits nicely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Footer Personal Information" text box (shape id 15)
$infoShape = $s.Shapes.Item(5)

# Correct the institute name - update the run text directly (not the
# paragraph) so the single run / formatting is preserved.
$infoShape.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Institute for Life Science & Technology"

# Reposition/resize the text box (EMU -> points, 12700 EMU per point).
# Add half an EMU before dividing so the float32 COM properties round-trip
# to the exact target EMU value instead of being truncated a hair short.
$infoShape.Left = (8610762 + 0.5) / 12700
$infoShape.Width = (7172497 + 0.5) / 12700

# "Footer Line - Logo/Info" connector (shape id 17)
$lineShape = $s.Shapes.Item(6)
$lineShape.Left = (8424252 + 0.5) / 12700
